# Update column G ("K", formerly "Strike#") values for rows 2-28
# with freshly regenerated strikeout counts, per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1, 0, 0, 2, 1, 2, 0, 1, 2, 1, 1, 1, 0, 0, 4, 2, 1, 1, 1, 1, 3, 2, 3, 5, 2, 3, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
